$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying weekly data rows (2, 3, 5) get rotated:
#   new row 2 <- old row 5
#   new row 3 <- old row 2
#   new row 5 <- old row 3
# Row 4 is untouched. Write the target values directly (captured from the
# diff) rather than shuffling live ranges, to avoid any read/overwrite
# ordering issues.

# Row 2 (was row 5's data)
$ws.Range("D2").Value = 44692
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

# Row 3 (was row 2's data)
$ws.Range("D3").Value = 44221
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1420
$ws.Range("N3").Value = "$/atado"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 1420
$ws.Range("Q3").Value = 1

# Row 5 (was row 3's data)
$ws.Range("D5").Value = 44691
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3250
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 542
$ws.Range("Q5").Value = 6
